$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("setback_distance")

# Insert a new log row at the top of the table (row 2), pushing existing
# entries down by one row.
$ws.Rows.Item(2).Insert()

# Copy date-formatting from the row below (now row 3) so the new date cell
# picks up the same number format style used throughout column A.
$ws.Range("A3").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A2").Value = 45180
$ws.Range("B2").Value = "Cell fraction"
$ws.Range("C2").Value = "When plotting the cell fraction I realized there were spots >1 which isn't physical. I traced it back to the inter quantile calculation which had an error where when the WSE was greater than the top WSE it was calculating a fraction much greater than 1. I corrected this and started the re-run.`n-> after update the recharge results were reduced on the high end so the flier on the maximum was outside the quartile. Also the effective recharge became less valuable at 1200m such that 0-600m are most valuable with 0,200 m gone then it is 600 followed by 1200 m."

# No content in column D for this new entry - remove the auto-carried format.
$ws.Range("D2").Clear()

$ws.Rows.Item(2).RowHeight = 120

# Update the active selection to match the new cursor position.
$ws.Range("C3").Select() | Out-Null
